$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '250.21'
Set-TextCell 'G2' '2'
Set-TextCell 'D3' '22.69'
Set-TextCell 'G3' '2'
Set-TextCell 'D4' '5.451'
Set-TextCell 'G4' '2'
Set-TextCell 'D5' '0.05728'
Set-TextCell 'G5' '2'
Set-TextCell 'D6' '3.412'
Set-TextCell 'G6' '2'
Set-TextCell 'D7' '6.334'
Set-TextCell 'G7' '2'
Set-TextCell 'D8' '0.8130'
Set-TextCell 'G8' '2'
Set-TextCell 'D9' '0.9331'
Set-TextCell 'G9' '2'
Set-TextCell 'D10' '0.1423'
Set-TextCell 'G10' '2'
Set-TextCell 'D11' '0.07526'
Set-TextCell 'G11' '2'
Set-TextCell 'G12' '2'
Set-TextCell 'D13' '0.03089'
Set-TextCell 'G13' '2'
Set-TextCell 'D14' '0.09366'
Set-TextCell 'G14' '2'
Set-TextCell 'D15' '3.722'
Set-TextCell 'G15' '2'
Set-TextCell 'D16' '0.001578'
Set-TextCell 'G16' '2'
Set-TextCell 'G17' '2'
Set-TextCell 'D18' '0.0005793'
Set-TextCell 'G18' '2'
Set-TextCell 'D19' '0.006435'
Set-TextCell 'G19' '2'
Set-TextCell 'D20' '0.005017'
Set-TextCell 'G20' '2'
Set-TextCell 'D21' '0.001026'
Set-TextCell 'G21' '2'
Set-TextCell 'D22' '0.0001500'
Set-TextCell 'G22' '2'
Set-TextCell 'D23' '3.701'
Set-TextCell 'G23' '2'
Set-TextCell 'D24' '2.165'
Set-TextCell 'G24' '2'
Set-TextCell 'G25' '2'
Set-TextCell 'D26' '0.1308'
Set-TextCell 'G26' '2'
Set-TextCell 'G27' '2'
Set-TextCell 'D28' '0.0003002'
Set-TextCell 'G28' '2'
Set-TextCell 'G29' '2'
Set-TextCell 'G30' '2'
Set-TextCell 'G31' '2'
Set-TextCell 'G32' '2'
Set-TextCell 'G33' '2'
Set-TextCell 'G34' '2'
Set-TextCell 'G35' '2'
Set-TextCell 'G36' '2'
Set-TextCell 'G37' '2'
Set-TextCell 'G38' '2'
Set-TextCell 'G39' '2'
Set-TextCell 'D40' '0.04009'
Set-TextCell 'G40' '2'
Set-TextCell 'B41' 'KickToken'
Set-TextCell 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell 'D41' '0.006872'
Set-TextCell 'E41' '40KickTokenKICK'
Set-TextCell 'G41' '2'
Set-TextCell 'B42' 'BKEXToken'
Set-TextCell 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell 'D42' '0.1068'
Set-TextCell 'E42' '41BKEXTokenBKK'
Set-TextCell 'G42' '2'
Set-TextCell 'B43' 'CEJI'
Set-TextCell 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextCell 'D43' '0.002711'
Set-TextCell 'E43' '42CEJICEJI'
Set-TextCell 'G43' '2'
Set-TextCell 'D44' '0.007967'
Set-TextCell 'G44' '2'
Set-TextCell 'D45' '0.00005893'
Set-TextCell 'G45' '2'
Set-TextCell 'D46' '0.00000000750'
Set-TextCell 'G46' '2'
Set-TextCell 'D47' '0.5003'
Set-TextCell 'E47' '46CoinbaseStockTokenCOINWorstin24h'
Set-TextCell 'G47' '2'
Set-TextCell 'G48' '2'
Set-TextCell 'D49' '0.00002101'
Set-TextCell 'G49' '2'
Set-TextCell 'D50' '0.01011'
Set-TextCell 'G50' '2'
Set-TextCell 'G51' '2'
